$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 13.23746348673675
$ws.Range("D2").Value = 7.861840526418574
$ws.Range("E2").Value = 13.74095183263102
$ws.Range("F2").Value = 41.05149137349137
$ws.Range("G2").Value = 3.713539489774528
$ws.Range("J2").Value = 10.88061863655646
$ws.Range("K2").Value = 17.75223641474558
$ws.Range("L2").Value = 9.498150120840561
$ws.Range("M2").Value = 19.82269561601675
$ws.Range("O2").Value = 31.77028743180739
$ws.Range("C3").Value = 13.22425584661412
$ws.Range("D3").Value = 7.838529891602969
$ws.Range("E3").Value = 13.76651090745077
$ws.Range("F3").Value = 41.21954515935457
$ws.Range("G3").Value = 3.715969857776541
$ws.Range("J3").Value = 10.91018932197322
$ws.Range("K3").Value = 17.33209219143146
$ws.Range("L3").Value = 9.513395207208077
$ws.Range("M3").Value = 19.65887455561484
$ws.Range("O3").Value = 31.92225163239497
$ws.Range("C4").Value = 13.21874699075762
$ws.Range("D4").Value = 7.824916626170099
$ws.Range("E4").Value = 13.78393265460808
$ws.Range("F4").Value = 41.33264472467492
$ws.Range("G4").Value = 3.717540565344657
$ws.Range("J4").Value = 10.92938439290766
$ws.Range("K4").Value = 17.0702152158143
$ws.Range("L4").Value = 9.523297523306272
$ws.Range("M4").Value = 19.55966481800317
$ws.Range("O4").Value = 32.02281814851498
$ws.Range("C5").Value = 13.2171583315758
$ws.Range("D5").Value = 7.819548088194829
$ws.Range("E5").Value = 13.79146709523165
$ws.Range("F5").Value = 41.38121992805446
$ws.Range("G5").Value = 3.718200433273164
$ws.Range("J5").Value = 10.93746831961207
$ws.Range("K5").Value = 16.9626546028017
$ws.Range("L5").Value = 9.527469423672184
$ws.Range("M5").Value = 19.51961410422387
$ws.Range("O5").Value = 32.06562202797725
$ws.Range("C6").Value = 13.21693422983201
$ws.Range("D6").Value = 7.818667543464977
$ws.Range("E6").Value = 13.79274446158474
$ws.Range("F6").Value = 41.38943576654809
$ws.Range("G6").Value = 3.718311201115412
$ws.Range("J6").Value = 10.93882647881477
$ws.Range("K6").Value = 16.94474747631206
$ws.Range("L6").Value = 9.528170428009968
$ws.Range("M6").Value = 19.51298746367785
$ws.Range("O6").Value = 32.07283953935064
$ws.Range("C7").Value = 13.21872290573182
$ws.Range("D7").Value = 7.824843495549351
$ws.Range("E7").Value = 13.78403250524959
$ws.Range("F7").Value = 41.33328976927346
$ws.Range("G7").Value = 3.71754938433225
$ws.Range("J7").Value = 10.92949235480374
$ws.Range("K7").Value = 17.06876784337999
$ws.Range("L7").Value = 9.52335323329676
$ws.Range("M7").Value = 19.55912310782578
$ws.Range("O7").Value = 32.02338804299479
$ws.Range("C8").Value = 13.23237120955236
$ws.Range("D8").Value = 7.853659703423562
$ws.Range("E8").Value = 13.74940607834237
$ws.Range("F8").Value = 41.10737469779865
$ws.Range("G8").Value = 3.714361236896412
$ws.Range("J8").Value = 10.89059943578238
$ws.Range("K8").Value = 17.6082626795191
$ws.Range("L8").Value = 9.503294418865302
$ws.Range("M8").Value = 19.76594355852525
$ws.Range("O8").Value = 31.82117588061772
$ws.Range("C9").Value = 13.27965371564855
$ws.Range("D9").Value = 7.915579682869946
$ws.Range("E9").Value = 13.69520307080965
$ws.Range("F9").Value = 40.74331366079845
$ws.Range("G9").Value = 3.708728822566887
$ws.Range("J9").Value = 10.82254237980973
$ws.Range("K9").Value = 18.62947220592374
$ws.Range("L9").Value = 9.468240014639356
$ws.Range("M9").Value = 20.18102388350644
$ws.Range("O9").Value = 31.48239271651254
$ws.Range("C10").Value = 13.32671905572647
$ws.Range("D10").Value = 7.964178591637047
$ws.Range("E10").Value = 13.66371153007691
$ws.Range("F10").Value = 40.52437071367693
$ws.Range("G10").Value = 3.704964242228202
$ws.Range("J10").Value = 10.77750706140451
$ws.Range("K10").Value = 19.35020963217931
$ws.Range("L10").Value = 9.445070450020193
$ws.Range("M10").Value = 20.48985184724731
$ws.Range("O10").Value = 31.26890417293001
$ws.Range("C11").Value = 13.35075948485443
$ws.Range("D11").Value = 7.986918580165485
$ws.Range("E11").Value = 13.65118961650668
$ws.Range("F11").Value = 40.43538911974638
$ws.Range("G11").Value = 3.703331871390429
$ws.Range("J11").Value = 10.75808924591114
$ws.Range("K11").Value = 19.67032502344515
$ws.Range("L11").Value = 9.435086035842161
$ws.Range("M11").Value = 20.63077162544942
$ws.Range("O11").Value = 31.17951832183011
$ws.Range("C12").Value = 13.36023644591956
$ws.Range("D12").Value = 7.995616637768683
$ws.Range("E12").Value = 13.64670683980216
$ws.Range("F12").Value = 40.40322743183099
$ws.Range("G12").Value = 3.702725194536896
$ws.Range("J12").Value = 10.7508893124155
$ws.Range("K12").Value = 19.79032816467321
$ws.Range("L12").Value = 9.431384682457329
$ws.Range("M12").Value = 20.68416157470433
$ws.Range("O12").Value = 31.14678593822322
$ws.Range("C13").Value = 13.35817888850253
$ws.Range("D13").Value = 7.993739554149803
$ws.Range("E13").Value = 13.64766077316711
$ws.Range("F13").Value = 40.41008570422454
$ws.Range("G13").Value = 3.702855344220582
$ws.Range("J13").Value = 10.75243314317464
$ws.Range("K13").Value = 19.76453904545742
$ws.Range("L13").Value = 9.432178303792918
$ws.Range("M13").Value = 20.67266246957095
$ws.Range("O13").Value = 31.15378573390558
$ws.Range("C14").Value = 13.35153170642861
$ws.Range("D14").Value = 7.987632455779931
$ws.Range("E14").Value = 13.65081562727621
$ws.Range("F14").Value = 40.4327123807393
$ws.Range("G14").Value = 3.703281730254565
$ws.Range("J14").Value = 10.7574938363592
$ws.Range("K14").Value = 19.68022271693208
$ws.Range("L14").Value = 9.43477993146014
$ws.Range("M14").Value = 20.63516371633002
$ws.Range("O14").Value = 31.1768030154911
$ws.Range("C15").Value = 13.3475085887154
$ws.Range("D15").Value = 7.983902881155587
$ws.Range("E15").Value = 13.65278178432957
$ws.Range("F15").Value = 40.44677181400191
$ws.Range("G15").Value = 3.703544395530087
$ws.Range("J15").Value = 10.76061358847332
$ws.Range("K15").Value = 19.62841504907851
$ws.Range("L15").Value = 9.436383849791829
$ws.Range("M15").Value = 20.61219707317379
$ws.Range("O15").Value = 31.19104724436988
$ws.Range("C16").Value = 13.32520047847683
$ws.Range("D16").Value = 7.962704861659983
$ws.Range("E16").Value = 13.66456613104556
$ws.Range("F16").Value = 40.53040011083154
$ws.Range("G16").Value = 3.705072529322444
$ws.Range("J16").Value = 10.77879752504566
$ws.Range("K16").Value = 19.32912436674446
$ws.Range("L16").Value = 9.445734103123892
$ws.Range("M16").Value = 20.48064807085723
$ws.Range("O16").Value = 31.2749016914824
$ws.Range("C17").Value = 13.31218551283213
$ws.Range("D17").Value = 7.949859661648697
$ws.Range("E17").Value = 13.67225716999863
$ws.Range("F17").Value = 40.58442761176494
$ws.Range("G17").Value = 3.706030477059404
$ws.Range("J17").Value = 10.7902261812086
$ws.Range("K17").Value = 19.14345796294389
$ws.Range("L17").Value = 9.451612211384909
$ws.Range("M17").Value = 20.40003245589204
$ws.Range("O17").Value = 31.32832690696729
$ws.Range("C18").Value = 13.30494752966161
$ws.Range("D18").Value = 7.942531183411705
$ws.Range("E18").Value = 13.67685065360552
$ws.Range("F18").Value = 40.61650167963528
$ws.Range("G18").Value = 3.706589011432655
$ws.Range("J18").Value = 10.79690029127275
$ws.Range("K18").Value = 19.03594355632607
$ws.Range("L18").Value = 9.455045454764567
$ws.Range("M18").Value = 20.35370730033599
$ws.Range("O18").Value = 31.35978310686314
$ws.Range("C19").Value = 13.30253958737254
$ws.Range("D19").Value = 7.940060266720803
$ws.Range("E19").Value = 13.67843510278755
$ws.Range("F19").Value = 40.62753274866003
$ws.Range("G19").Value = 3.706779419807012
$ws.Range("J19").Value = 10.79917733363094
$ws.Range("K19").Value = 18.99941994997396
$ws.Range("L19").Value = 9.456216887400348
$ws.Range("M19").Value = 20.33803083642987
$ws.Range("O19").Value = 31.37055844001462
$ws.Range("C20").Value = 13.31354535716495
$ws.Range("D20").Value = 7.951220900647606
$ws.Range("E20").Value = 13.67142087459927
$ws.Range("F20").Value = 40.57857286022941
$ws.Range("G20").Value = 3.705927721113764
$ws.Range("J20").Value = 10.78899916839397
$ws.Range("K20").Value = 19.16329817196095
$ws.Range("L20").Value = 9.450981065337634
$ws.Range("M20").Value = 20.40860995235046
$ws.Range("O20").Value = 31.32256438386155
$ws.Range("C21").Value = 13.35347405238037
$ws.Range("D21").Value = 7.989423931887801
$ws.Range("E21").Value = 13.64988194409975
$ws.Range("F21").Value = 40.42602469760557
$ws.Range("G21").Value = 3.703156179615556
$ws.Range("J21").Value = 10.75600323537331
$ws.Range("K21").Value = 19.70502229127865
$ws.Range("L21").Value = 9.434013614832539
$ws.Range("M21").Value = 20.64617756656615
$ws.Range("O21").Value = 31.17001195633192
$ws.Range("C22").Value = 13.38174352425357
$ws.Range("D22").Value = 8.014896529554083
$ws.Range("E22").Value = 13.63731446435597
$ws.Range("F22").Value = 40.335268435636
$ws.Range("G22").Value = 3.701411625467829
$ws.Range("J22").Value = 10.73533103999912
$ws.Range("K22").Value = 20.05193122043998
$ws.Range("L22").Value = 9.423387788643392
$ws.Range("M22").Value = 20.80158101626889
$ws.Range("O22").Value = 31.07681765785814
$ws.Range("C23").Value = 13.36645838282143
$ws.Range("D23").Value = 8.001256483241109
$ws.Range("E23").Value = 13.64388398119293
$ws.Range("F23").Value = 40.38288621342386
$ws.Range("G23").Value = 3.7023366334992
$ws.Range("J23").Value = 10.74628268943863
$ws.Range("K23").Value = 19.86746449602748
$ws.Range("L23").Value = 9.429016709930185
$ws.Range("M23").Value = 20.71863818646298
$ws.Range("O23").Value = 31.12596031527751
$ws.Range("C24").Value = 13.31292980961965
$ws.Range("D24").Value = 7.950605308539053
$ws.Range("E24").Value = 13.67179842857285
$ws.Range("F24").Value = 40.58121663861668
$ws.Range("G24").Value = 3.705974152783283
$ws.Range("J24").Value = 10.78955357805862
$ws.Range("K24").Value = 19.15433081115388
$ws.Range("L24").Value = 9.451266238727895
$ws.Range("M24").Value = 20.40473199546499
$ws.Range("O24").Value = 31.32516731160754
$ws.Range("C25").Value = 13.2646829527785
$ws.Range("D25").Value = 7.898268749243802
$ws.Range("E25").Value = 13.70840178907908
$ws.Range("F25").Value = 40.83330724791555
$ws.Range("G25").Value = 3.710186643014869
$ws.Range("J25").Value = 10.84007863863504
$ws.Range("K25").Value = 18.35791016037834
$ws.Range("L25").Value = 9.477267440629793
$ws.Range("M25").Value = 20.06791375764936
$ws.Range("O25").Value = 31.56783837514572
